$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old contents of A1/A3 from the previous version.
$ws.Range("A1").ClearContents()
$ws.Range("A3").ClearContents()

# Row 1: a new "deleted" marker in column B.
$ws.Range("B1").Value = "deleted"

# Row 3: earlier cell's text was modified, with a "Modified" marker in column B.
$ws.Range("A3").Value = "#3. this earlier cell was modified to this text"
$ws.Range("B3").Value = "Modified"

# Row 5: a new cell added in the second version, with an "Added" marker in column B.
$ws.Range("A5").Value = "#2. this cell was added in the second version"
$ws.Range("B5").Value = "Added"

# Update the active selection to match the new version.
$ws.Range("E15").Select()
